# Updated Gantt Chart dates and who's implementing
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Credit the teammate implementing each feature (append " - Name") ---
$ws.Cells.Item(5, 2).Value = "Login - Nhan"
$ws.Cells.Item(6, 2).Value = "Logout - Nhan"
$ws.Cells.Item(7, 2).Value = "Create New Account - Nhan"
$ws.Cells.Item(8, 2).Value = "Delete Account - Nhan"
$ws.Cells.Item(9, 2).Value = "Update Account - Nhan"
$ws.Cells.Item(10, 2).Value = "Forgot Password - Nhan"
$ws.Cells.Item(11, 2).Value = "Seller Create Item - Aaron"
$ws.Cells.Item(12, 2).Value = "Seller Update Item - Aaron"
$ws.Cells.Item(13, 2).Value = "Seller Delete Item - Aaron"
$ws.Cells.Item(14, 2).Value = "View All Items - Aaron"
$ws.Cells.Item(15, 2).Value = "View Single Item - Aaron"

# --- Rows 16-22: re-ordered tasks, refreshed dates/% complete, new owners ---

# Row 16: Add Item to Cart - Anh
$ws.Cells.Item(16, 2).Value = "Add Item to Cart - Anh"
$ws.Cells.Item(16, 3).Value = 4
$ws.Cells.Item(16, 4).Value = 13
$ws.Cells.Item(16, 5).Value = 6
$ws.Cells.Item(16, 6).Value = 9
$ws.Cells.Item(16, 7).Value = 0.5

# Row 17: Remove Items from Cart - Anh
$ws.Cells.Item(17, 2).Value = "Remove Items from Cart - Anh"
$ws.Cells.Item(17, 3).Value = 10
$ws.Cells.Item(17, 4).Value = 13
$ws.Cells.Item(17, 5).ClearContents()
$ws.Cells.Item(17, 6).ClearContents()
$ws.Cells.Item(17, 7).Value = 0

# Row 18: Buy Items from Cart - Anh
$ws.Cells.Item(18, 2).Value = "Buy Items from Cart - Anh"
$ws.Cells.Item(18, 3).Value = 12
$ws.Cells.Item(18, 4).Value = 15
$ws.Cells.Item(18, 5).ClearContents()
$ws.Cells.Item(18, 6).ClearContents()
$ws.Cells.Item(18, 7).Value = 0

# Row 19: Maintain UI Consistency - Nhan
$ws.Cells.Item(19, 2).Value = "Maintain UI Consistency - Nhan"
$ws.Cells.Item(19, 3).Value = 1
$ws.Cells.Item(19, 4).Value = 13
$ws.Cells.Item(19, 5).Value = 1
$ws.Cells.Item(19, 6).Value = 13
$ws.Cells.Item(19, 7).Value = 0.6

# Row 20: Add Item Categories - Aaron
$ws.Cells.Item(20, 2).Value = "Add Item Categories - Aaron"
$ws.Cells.Item(20, 3).Value = 7
$ws.Cells.Item(20, 4).Value = 7
$ws.Cells.Item(20, 5).Value = 7
$ws.Cells.Item(20, 6).Value = 7
$ws.Cells.Item(20, 7).Value = 1

# Row 21: Add Pictures to Items - Anh
$ws.Cells.Item(21, 2).Value = "Add Pictures to Items - Anh"
$ws.Cells.Item(21, 3).Value = 11
$ws.Cells.Item(21, 4).Value = 16
$ws.Cells.Item(21, 5).ClearContents()
$ws.Cells.Item(21, 6).ClearContents()
$ws.Cells.Item(21, 7).Value = 0

# Row 22: View Purchase History - Aaron
$ws.Cells.Item(22, 2).Value = "View Purchase History - Aaron"
$ws.Cells.Item(22, 3).Value = 12
$ws.Cells.Item(22, 4).Value = 15
$ws.Cells.Item(22, 5).ClearContents()
$ws.Cells.Item(22, 6).ClearContents()
$ws.Cells.Item(22, 7).Value = 0

# --- Widen the Feature column to fit the longer labels ---
$ws.Columns.Item(2).ColumnWidth = 26.15

# --- Move the active selection ---
$ws.Range("E33").Select() | Out-Null
